# Analyducks "Ducks" sheet update:
#  - clear the stray placeholder Name values (a/b/c/d/e) in B100:B104
#  - append 5 new duck purchases (rows 105-109: Flavortown + the Las Vegas trip)
#  - grow Table1 to cover the new rows
#  - nudge column K's width and the saved selection/scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ducks")

# 1. Remove the leftover "a","b","c","d","e" placeholder values from column B (Name)
#    for the last five pre-existing rows.
$ws.Range("B100:B104").ClearContents()

# 2. Stage formatting (number formats / styles) for the five new rows by copying
#    the style of the last existing row (104) down, row by row. Column B is
#    intentionally skipped so no cell is created there at all.
for ($r = 105; $r -le 109; $r++) {
    $ws.Range("A104").Copy($ws.Range("A$r"))
    $ws.Range("C104:R104").Copy($ws.Range("C$r`:R$r"))
}

# 3. Row 105 - Flavortown (Pigeon Forge, TN)
$ws.Range("A105").Value = "Flavortown"
$ws.Range("D105").Value = "Flavortown"
$ws.Range("E105").Value = "Pigeon Forge"
$ws.Range("F105").Value = "TN"
$ws.Range("I105").Value = 45197
$ws.Range("J105").Value = 35.823491320947902
$ws.Range("K105").Value = -83.562385482795094
$ws.Range("M105").Value = "Allan"
$ws.Range("N105").Value = 1
$ws.Range("O105").Value = 20
$ws.Range("P105").Formula = "=2+5/16"
$ws.Range("Q105").Formula = "=1+11/16"
$ws.Range("R105").Formula = "=2+7/16"

# 4. Row 106 - Las Vegas Cards (Duck Donuts, Las Vegas, NV)
$ws.Range("A106").Value = "Las Vegas Cards"
$ws.Range("D106").Value = "Duck Donuts"
$ws.Range("E106").Value = "Las Vegas"
$ws.Range("F106").Value = "NV"
$ws.Range("I106").Value = 45201
$ws.Range("J106").Value = 36.114481645923
$ws.Range("K106").Value = -115.171873618758
$ws.Range("M106").Value = "Allan"
$ws.Range("N106").Value = 1
$ws.Range("O106").Value = 70
$ws.Range("P106").Formula = "=2.625"
$ws.Range("Q106").Value = 3
$ws.Range("R106").Value = 2.625

# 5. Row 107 - Clown (Duck Donuts, Las Vegas, NV)
$ws.Range("A107").Value = "Clown"
$ws.Range("D107").Value = "Duck Donuts"
$ws.Range("E107").Value = "Las Vegas"
$ws.Range("F107").Value = "NV"
$ws.Range("I107").Value = 45201
$ws.Range("J107").Value = 36.114481645923
$ws.Range("K107").Value = -115.171873618758
$ws.Range("M107").Value = "Allan"
$ws.Range("N107").Value = 1
$ws.Range("O107").Value = 17
$ws.Range("P107").Value = 2.25
$ws.Range("Q107").Formula = "=1+14/16"
$ws.Range("R107").Formula = "=1+14/16"

# 6. Row 108 - Koala (Duck Donuts, Las Vegas, NV)
$ws.Range("A108").Value = "Koala"
$ws.Range("D108").Value = "Duck Donuts"
$ws.Range("E108").Value = "Las Vegas"
$ws.Range("F108").Value = "NV"
$ws.Range("I108").Value = 45201
$ws.Range("J108").Value = 36.114481645923
$ws.Range("K108").Value = -115.171873618758
$ws.Range("M108").Value = "Allan"
$ws.Range("N108").Value = 1
$ws.Range("O108").Value = 15
$ws.Range("P108").Value = 2
$ws.Range("Q108").Formula = "=1+13/16"
$ws.Range("R108").Value = 1.875

# 7. Row 109 - Elvis (Las Vegas Airport, Las Vegas, Nv)
$ws.Range("A109").Value = "Elvis"
$ws.Range("D109").Value = "Las Vegas Airport"
$ws.Range("E109").Value = "Las Vegas"
$ws.Range("F109").Value = "Nv"
$ws.Range("I109").Value = 45204
$ws.Range("J109").Value = 36.083568911204402
$ws.Range("K109").Value = -115.149256886367
$ws.Range("M109").Value = "Allan"
$ws.Range("N109").Value = 1
$ws.Range("O109").Value = 54
$ws.Range("P109").Value = 3
$ws.Range("Q109").Value = 2.625
$ws.Range("R109").Value = 3.5

# 8. Grow Table1 so it covers the freshly added rows.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:R109"))

# 9. Column K (Longitude) widened slightly to fit the new, longer values.
$ws.Columns.Item(11).ColumnWidth = 11.5

# 10. Restore the frozen-pane view and leave the selection where the editor
#     finished (Q103), scrolled down a bit further than before.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 84
$win.ScrollColumn = 1
$ws.Range("Q103").Select()
